$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 470.25
$ws.Cells.Item(2, 9).Value = 466
$ws.Cells.Item(2, 11).Value = 466
$ws.Cells.Item(2, 13).Value = -353

$ws.Cells.Item(4, 8).Value = 106.6
$ws.Cells.Item(4, 9).Value = 108.25
$ws.Cells.Item(4, 11).Value = 108.25
$ws.Cells.Item(4, 13).Value = 5.75

$ws.Cells.Item(43, 8).Value = 0
$ws.Cells.Item(43, 9).Value = 0
$ws.Cells.Item(43, 10).Value = 0
$ws.Cells.Item(43, 11).Value = 0
$ws.Cells.Item(43, 12).Value = 0
$ws.Cells.Item(43, 13).Value = $null
$ws.Cells.Item(43, 14).Value = $null

$ws.Cells.Item(52, 8).Value = 1499
$ws.Cells.Item(52, 10).Value = 1499
$ws.Cells.Item(52, 12).Value = 4497
$ws.Cells.Item(52, 14).Value = -4817

$ws.Cells.Item(58, 8).Value = 564.5
$ws.Cells.Item(58, 10).Value = 989
$ws.Cells.Item(58, 12).Value = 2967
$ws.Cells.Item(58, 14).Value = -3267

$ws.Cells.Item(93, 8).Value = 0
$ws.Cells.Item(93, 10).Value = 0
$ws.Cells.Item(93, 12).Value = 0
$ws.Cells.Item(93, 14).Value = $null

$ws.Cells.Item(106, 8).Value = 3003
$ws.Cells.Item(106, 10).Value = 2006
$ws.Cells.Item(106, 12).Value = 2006
$ws.Cells.Item(106, 14).Value = -3268

$ws.Cells.Item(138, 8).Value = 4765267
$ws.Cells.Item(138, 9).Value = 16669457
$ws.Cells.Item(138, 10).Value = 3590.6
$ws.Cells.Item(138, 11).Value = 50008371
$ws.Cells.Item(138, 12).Value = 10771.8
$ws.Cells.Item(138, 13).Value = -50003231
$ws.Cells.Item(138, 14).Value = -21051.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(4, 8).Value = 499.42856
$ws.Cells.Item(4, 9).Value = 299.8
$ws.Cells.Item(4, 10).Value = 998.5
$ws.Cells.Item(4, 11).Value = 299.8
$ws.Cells.Item(4, 12).Value = 998.5
$ws.Cells.Item(4, 13).Value = -183.8
$ws.Cells.Item(4, 14).Value = -1230.5

$ws.Cells.Item(12, 8).Value = 1900
$ws.Cells.Item(12, 10).Value = 1900
$ws.Cells.Item(12, 12).Value = 1900
$ws.Cells.Item(12, 14).Value = -2246

$ws.Cells.Item(15, 8).Value = 3249.75
$ws.Cells.Item(15, 10).Value = 3249.75
$ws.Cells.Item(15, 12).Value = 3249.75
$ws.Cells.Item(15, 14).Value = -3949.75

$ws.Cells.Item(19, 8).Value = 4999
$ws.Cells.Item(19, 10).Value = 4999
$ws.Cells.Item(19, 12).Value = 4999
$ws.Cells.Item(19, 14).Value = -5457

$ws.Cells.Item(92, 8).Value = 0
$ws.Cells.Item(92, 10).Value = 0
$ws.Cells.Item(92, 12).Value = 0
$ws.Cells.Item(92, 14).Value = $null

$ws.Cells.Item(101, 8).Value = 29602
$ws.Cells.Item(101, 10).Value = 29602
$ws.Cells.Item(101, 12).Value = 29602
$ws.Cells.Item(101, 14).Value = -36092

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(22, 8).Value = 257.66666
$ws.Cells.Item(22, 9).Value = 257.66666
$ws.Cells.Item(22, 11).Value = 257.66666
$ws.Cells.Item(22, 13).Value = -84.66665999999998

$ws.Cells.Item(26, 8).Value = 2999
$ws.Cells.Item(26, 9).Value = 2999
$ws.Cells.Item(26, 11).Value = 2999
$ws.Cells.Item(26, 13).Value = -2707

$ws.Cells.Item(107, 8).Value = 974.5
$ws.Cells.Item(107, 9).Value = 974.5
$ws.Cells.Item(107, 10).Value = 0
$ws.Cells.Item(107, 11).Value = 974.5
$ws.Cells.Item(107, 12).Value = 0
$ws.Cells.Item(107, 13).Value = 945.5
$ws.Cells.Item(107, 14).Value = $null

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(4, 8).Value = 0
$ws.Cells.Item(4, 9).Value = 0
$ws.Cells.Item(4, 10).Value = 0
$ws.Cells.Item(4, 11).Value = 0
$ws.Cells.Item(4, 12).Value = 0
$ws.Cells.Item(4, 13).Value = $null
$ws.Cells.Item(4, 14).Value = $null

$ws.Cells.Item(19, 8).Value = 301.875
$ws.Cells.Item(19, 10).Value = 426.25
$ws.Cells.Item(19, 12).Value = 426.25
$ws.Cells.Item(19, 14).Value = -766.25

$ws.Cells.Item(22, 8).Value = 475
$ws.Cells.Item(22, 9).Value = 370
$ws.Cells.Item(22, 10).Value = 1000
$ws.Cells.Item(22, 11).Value = 370
$ws.Cells.Item(22, 12).Value = 1000
$ws.Cells.Item(22, 13).Value = -20
$ws.Cells.Item(22, 14).Value = -1700

$ws.Cells.Item(24, 8).Value = 301.875
$ws.Cells.Item(24, 10).Value = 426.25
$ws.Cells.Item(24, 12).Value = 426.25
$ws.Cells.Item(24, 14).Value = -766.25

$ws.Cells.Item(32, 8).Value = 4873.75
$ws.Cells.Item(32, 9).Value = 4500
$ws.Cells.Item(32, 10).Value = 4998.3335
$ws.Cells.Item(32, 11).Value = 4500
$ws.Cells.Item(32, 12).Value = 4998.3335
$ws.Cells.Item(32, 13).Value = -4184
$ws.Cells.Item(32, 14).Value = -5630.3335

$ws.Cells.Item(42, 8).Value = 8999
$ws.Cells.Item(42, 10).Value = 8999
$ws.Cells.Item(42, 12).Value = 8999
$ws.Cells.Item(42, 14).Value = -10185

$ws.Cells.Item(58, 8).Value = 6961.923
$ws.Cells.Item(58, 10).Value = 9332.333000000001
$ws.Cells.Item(58, 12).Value = 9332.333000000001
$ws.Cells.Item(58, 14).Value = -9738.333000000001

$ws.Cells.Item(122, 8).Value = 1936
$ws.Cells.Item(122, 9).Value = 1842.1666
$ws.Cells.Item(122, 11).Value = 5526.4998
$ws.Cells.Item(122, 13).Value = -3076.4998

$ws.Cells.Item(136, 8).Value = 6961.923
$ws.Cells.Item(136, 10).Value = 9332.333000000001
$ws.Cells.Item(136, 12).Value = 27996.999
$ws.Cells.Item(136, 14).Value = -33096.999

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 679
$ws.Cells.Item(2, 9).Value = 396.6154
$ws.Cells.Item(2, 10).Value = 1902.6666
$ws.Cells.Item(2, 11).Value = 2379.6924
$ws.Cells.Item(2, 12).Value = 11415.9996
$ws.Cells.Item(2, 13).Value = -2266.6924
$ws.Cells.Item(2, 14).Value = -11641.9996

$ws.Cells.Item(13, 8).Value = 5166.6665
$ws.Cells.Item(13, 9).Value = 4750
$ws.Cells.Item(13, 10).Value = 6000
$ws.Cells.Item(13, 11).Value = 14250
$ws.Cells.Item(13, 12).Value = 18000
$ws.Cells.Item(13, 13).Value = -14082
$ws.Cells.Item(13, 14).Value = -18336

$ws.Cells.Item(17, 8).Value = 620
$ws.Cells.Item(17, 10).Value = 900
$ws.Cells.Item(17, 12).Value = 2700
$ws.Cells.Item(17, 14).Value = -3038

$ws.Cells.Item(141, 8).Value = 0
$ws.Cells.Item(141, 9).Value = 0
$ws.Cells.Item(141, 11).Value = 0
$ws.Cells.Item(141, 13).Value = $null

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(5, 8).Value = 4
$ws.Cells.Item(5, 9).Value = 4
$ws.Cells.Item(5, 11).Value = 4
$ws.Cells.Item(5, 13).Value = 108

$ws.Cells.Item(80, 8).Value = 2666.3333
$ws.Cells.Item(80, 9).Value = 0
$ws.Cells.Item(80, 10).Value = 2666.3333
$ws.Cells.Item(80, 11).Value = 0
$ws.Cells.Item(80, 12).Value = 2666.3333
$ws.Cells.Item(80, 13).Value = $null
$ws.Cells.Item(80, 14).Value = -4662.3333

$ws.Cells.Item(83, 8).Value = 2666.3333
$ws.Cells.Item(83, 9).Value = 0
$ws.Cells.Item(83, 10).Value = 2666.3333
$ws.Cells.Item(83, 11).Value = 0
$ws.Cells.Item(83, 12).Value = 13331.6665
$ws.Cells.Item(83, 13).Value = $null
$ws.Cells.Item(83, 14).Value = -23315.6665

$ws.Cells.Item(133, 8).Value = 95000
$ws.Cells.Item(133, 9).Value = 95000
$ws.Cells.Item(133, 11).Value = 95000
$ws.Cells.Item(133, 13).Value = -89940

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(2, 8).Value = 0
$ws.Cells.Item(2, 10).Value = 0
$ws.Cells.Item(2, 12).Value = 0
$ws.Cells.Item(2, 14).Value = $null

$ws.Cells.Item(100, 8).Value = 0
$ws.Cells.Item(100, 9).Value = 0
$ws.Cells.Item(100, 11).Value = 0
$ws.Cells.Item(100, 13).Value = $null

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(2, 8).Value = 12500
